$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 8.937933333333332
$ws.Range("H2").Value = 26.8138
$ws.Range("I2").Value = 0.2302024600837126
$ws.Range("J2").Value = 0.2302024600837126
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723262999999999
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 28.96862549215555
$ws.Range("R2").Value = 260.7176294293999
$ws.Range("S2").Value = 0.005723426797525321
$ws.Range("T2").Value = 0.005723426797525321
$ws.Range("G3").Value = 8.937933333333332
$ws.Range("H3").Value = 26.8138
$ws.Range("I3").Value = 0.2302024600837126
$ws.Range("J3").Value = 0.2302024600837126
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 24.42298327473333
$ws.Range("R3").Value = 219.8068494726
$ws.Range("S3").Value = 0.004825329285573923
$ws.Range("T3").Value = 0.004825329285573922
$ws.Range("G4").Value = 8.937933333333332
$ws.Range("H4").Value = 26.8138
$ws.Range("I4").Value = 0.2302024600837126
$ws.Range("J4").Value = 0.2302024600837126
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 1111.758062828689
$ws.Range("R4").Value = 10005.8225654582
$ws.Range("S4").Value = 0.2196537040006133
$ws.Range("T4").Value = 0.2196537040006133
$ws.Range("I5").Value = 0.5278886986241245
$ws.Range("J5").Value = 0.5278886986241244
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723262999999999
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 66.42939439666567
$ws.Range("R5").Value = 597.8645495699909
$ws.Range("S5").Value = 0.013124674352817
$ws.Range("T5").Value = 0.013124674352817
$ws.Range("I6").Value = 0.5278886986241245
$ws.Range("J6").Value = 0.5278886986241244
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("S6").Value = 0.0110652023269786
$ws.Range("T6").Value = 0.01106520232697859
$ws.Range("I7").Value = 0.5278886986241245
$ws.Range("J7").Value = 0.5278886986241244
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 2549.427650591114
$ws.Range("R7").Value = 22944.84885532002
$ws.Range("S7").Value = 0.5036988219443288
$ws.Range("T7").Value = 0.5036988219443288
$ws.Range("G8").Value = 9.392449999999998
$ws.Range("H8").Value = 28.17735
$ws.Range("I8").Value = 0.241908841292163
$ws.Range("J8").Value = 0.2419088412921629
$ws.Range("M8").Value = 3.241087666666667
$ws.Range("N8").Value = 9.723262999999999
$ws.Range("O8").Value = 0.02486257877280725
$ws.Range("P8").Value = 0.02486257877280725
$ws.Range("Q8").Value = 30.44175385478333
$ws.Range("R8").Value = 273.9757846930499
$ws.Range("S8").Value = 0.006014477622464928
$ws.Range("T8").Value = 0.006014477622464929
$ws.Range("G9").Value = 9.392449999999998
$ws.Range("H9").Value = 28.17735
$ws.Range("I9").Value = 0.241908841292163
$ws.Range("J9").Value = 0.2419088412921629
$ws.Range("O9").Value = 0.02096124117795788
$ws.Range("P9").Value = 0.02096124117795788
$ws.Range("Q9").Value = 25.66495415704999
$ws.Range("R9").Value = 230.98458741345
$ws.Range("S9").Value = 0.005070709565405365
$ws.Range("T9").Value = 0.005070709565405365
$ws.Range("G10").Value = 9.392449999999998
$ws.Range("H10").Value = 28.17735
$ws.Range("I10").Value = 0.241908841292163
$ws.Range("J10").Value = 0.2419088412921629
$ws.Range("M10").Value = 124.3864796666667
$ws.Range("N10").Value = 373.159439
$ws.Range("O10").Value = 0.9541761800492348
$ws.Range("P10").Value = 0.9541761800492349
$ws.Range("Q10").Value = 1168.293790945183
$ws.Range("R10").Value = 10514.64411850665
$ws.Range("S10").Value = 0.2308236541042926
$ws.Range("T10").Value = 0.2308236541042926
